$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.09889220731696108, 2.044145585884516, 20.46097368034341, 4.523380779941415, 4.62393703887486, 23)
    3  = @(1.536575381265211, 2.134447601952433, 22.16480059406862, 4.707950785009187, 4.554862431766541, 22)
    4  = @(0.2287529840863657, 1.32172702803788, 4.6367541873095, 2.153312375692273, 2.194002637318409, 21)
    5  = @(0.7295427155613371, 0.9049081462973387, 1.224543713265862, 1.106591032525504, 0.8536677151144266, 20)
    6  = @(0.6468487310729173, 0.8752534674755443, 1.269018386393639, 1.126507162158164, 0.9475552462509015, 19)
    7  = @(0.2186385474789455, 0.6055429715130134, 0.4952421196251417, 0.7037344098629409, 0.6883017247615748, 18)
    8  = @(0.198091363434846, 0.5302076377283306, 0.4462351981637913, 0.6680083818065394, 0.6575957709830825, 17)
    9  = @(0.2488500629165574, 0.2488500629165574, 0.0794481986358422, 0.2818655683758522, 0.1388309378506621, 11)
    10 = @(0.2825386127611949, 0.2825386127611949, 0.1014758371739403, 0.318552722753927, 0.1589205180325681, 7)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
